$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 264
$ws.Range("C3").Value = 162374
$ws.Range("C4").Value = 153382
$ws.Range("C8").Value = 64.55
